$d = $word.ActiveDocument

# The final answer paragraph under "Output" ("Ingen kassasjoner er
# registrert.") carries a leftover _GoBack bookmark from the last edit
# position. Remove it - the new revision no longer marks that spot.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Remember which paragraph is the last one before we start appending -
# that's the "Ingen kassasjoner er registrert." paragraph.
$targetIndex = $d.Paragraphs.Count

# Append four new paragraphs after it:
#   (empty)
#   AND/OR
#   (empty)
#   Kassasjoner er registrert.
$idx = $targetIndex
for ($n = 0; $n -lt 4; $n++) {
    $p = $d.Paragraphs.Item($idx)
    $r = $p.Range
    $r.Collapse(0)            # wdCollapseEnd
    $r.InsertParagraphAfter()
    $idx = $idx + 1
}

# Fill in the text of the two non-empty new paragraphs.
$pAndOr = $d.Paragraphs.Item($targetIndex + 2)
$pAndOr.Range.InsertBefore("AND/OR")

$pFinal = $d.Paragraphs.Item($targetIndex + 4)
$pFinal.Range.InsertBefore("Kassasjoner er registrert.")
